$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copyright year bump ---
$ws.Range("B3").Value = "Copyright @2015 - 2023"

# --- Quarter header row (row 8): shift one quarter forward, add new quarter in column H ---
$ws.Range("D8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("E8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("F8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("G8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("H8").Value = "فصل سوم منتهی به 1401/09"

# --- Publish date row (row 9) ---
$ws.Range("D9").Value = "1401-11-02 (2)"
$ws.Range("E9").Value = "1401-11-02 (7)"
$ws.Range("F9").Value = "1401-04-29"
$ws.Range("G9").Value = "1401-09-16 (3)"
$ws.Range("H9").Value = "1401-11-02"

# --- Column widths: the 'current quarter' highlighted column moves from F to E ---
$ws.Columns("E").ColumnWidth = 30.166666666666668
$ws.Columns("F").ColumnWidth = 28.166666666666668

# --- Quarterly figures (D:H), each quarter shifts one column left, columns G/H carry new data ---
$rows = @{
    12 = @(14617,614021,978224,164105,1414170)
    13 = @(-86788,-75676,-92891,-9363,-27781)
    14 = @(-72171,538345,885333,154742,1386389)
    16 = @(0,5632,0,0,0)
    17 = @(-146578,-498083,-49963,-96917,-104625)
    18 = @(0,0,0,0,0)
    19 = @(0,0,0,0,0)
    20 = @(0,-17,-780,-1,455)
    21 = @(0,0,0,0,0)
    22 = @(0,0,0,0,0)
    23 = @(0,0,0,0,0)
    24 = @(0,0,-1878521,-58015,58016)
    25 = @(0,0,0,0,14751)
    26 = @(0,0,0,-4134,4134)
    27 = @(0,0,0,0,0)
    28 = @(0,0,0,0,0)
    29 = @(0,0,0,0,0)
    30 = @(432247,-432247,0,0,0)
    31 = @(8630,475254,14637,28820,97930)
    32 = @(294299,-449461,-1914627,-130247,70661)
    33 = @(222128,88884,-1029294,24495,1457050)
    35 = @(-3739,18729,0,0,0)
    36 = @("-","-","-","-","-")
    37 = @(0,0,0,0,0)
    38 = @(0,0,0,0,0)
    39 = @(3043165,5598948,4521700,4079099,3538576)
    40 = @(-3013384,-4299862,-4395469,-3284009,-3588265)
    41 = @(-308459,-340121,-548059,-688545,-510095)
    42 = @(0,2000000,0,0,0)
    43 = @(0,0,0,0,0)
    44 = @(0,0,0,0,0)
    45 = @(0,0,0,0,0)
    46 = @(0,0,0,0,0)
    47 = @(0,0,0,0,0)
    48 = @(0,0,0,0,0)
    49 = @(0,0,0,0,0)
    50 = @(-37358,-1675308,-1628,-261995,-503669)
    51 = @(-319775,1302386,-423456,-155450,-1063453)
    52 = @(-97647,1391270,-1452750,-130955,393597)
    53 = @(573634,475987,1868143,415393,283975)
    54 = @(0,886,0,-463,463)
    55 = @(475987,1868143,415393,283975,678035)
    56 = @(0,0,0,0,0)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $ws.Cells.Item([int]$r, 4 + $i).Value = $vals[$i]
    }
}
